$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 "Marking": Right count 5 -> 4, Wrong mark -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total": Right total 75 -> 60, Wrong total -6 -> -12, and the
# "score / max" summary text updated to match
$ws.Range("B12").Value = 60
$ws.Range("C12").Value = -12
$ws.Range("E12").Value = "48 / 112"
